$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record was inserted at row 58, pushing the existing
# rows 58-76 down to 59-77 (all their data stays the same, just shifted).
$ws.Rows.Item(58).EntireRow.Insert()

# Populate the newly inserted row 58 with the new record's data.
$ws.Range("A58").Value = 1
$ws.Range("B58").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C58").Value = "Arica y Parinacota"
$ws.Range("D58").Value = 45204
$ws.Range("E58").Value = 15
$ws.Range("F58").Value = 100112052
$ws.Range("G58").Value = "Albahaca"
$ws.Range("H58").Value = "Sin especificar"
$ws.Range("I58").Value = "Primera"
$ws.Range("J58").Value = 350
$ws.Range("K58").Value = 800
$ws.Range("L58").Value = 1000
$ws.Range("M58").Value = 914
$ws.Range("N58").Value = "$/paquete"
$ws.Range("O58").Value = "Región de Arica y Parinacota"
$ws.Range("P58").Value = 914
$ws.Range("Q58").Value = 1
$ws.Range("R58").Value = "Hortaliza"
